# ------------------------------------------------------------------
# Helper: insert a brand-new single-run paragraph right after $anchor
# (a Range) and return a Range over the paragraph that was inserted,
# so it can be used as the next anchor.
# ------------------------------------------------------------------
$d = $word.ActiveDocument

function New-ParaAfter($anchorRange, [string]$text) {
    $anchorRange.InsertParagraphAfter()
    $newPara = $d.Range($anchorRange.End + 1, $anchorRange.End + 1)
    $newPara.Text = $text
    return $newPara
}

# Append one more run of text right at the end of $paraRange (which must
# be the Range of the paragraph text inserted by New-ParaAfter / itself),
# forcing a *separate* <w:r> element (not merged with the previous run)
# by briefly wrapping the insertion point in a bookmark.
function Add-Run($paraRange, [string]$text) {
    $mark = $d.Bookmarks.Add("__tmp_split__", $d.Range($paraRange.End, $paraRange.End))
    $ins = $d.Range($paraRange.End, $paraRange.End)
    $ins.InsertAfter($text)
    $d.Bookmarks("__tmp_split__").Delete()
    return $d.Range($paraRange.Start, $paraRange.End + $text.Length)
}

# ------------------------------------------------------------------
# 1) "- Dùng để xem tour du lịch, đặt tour du lịch" after paragraph 1
# ------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("1. Website dùng để làm gì?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = New-ParaAfter $rng "- Dùng để xem tour du lịch, đặt tour du lịch"

# ------------------------------------------------------------------
# 2) "- Khách hàng trung lưu đi theo nhóm 1 vài gia đình" after paragraph 2
# ------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("2. Đối tượng sử dụng website là ai?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$null = New-ParaAfter $rng "- Khách hàng trung lưu đi theo nhóm 1 vài gia đình"

# ------------------------------------------------------------------
# 3) Four new paragraphs after paragraph "3. Website có những chức năng gì?"
# ------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("3. Website có những chức năng gì?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# 3a) "- Xem các tour " / "sắp có" / "."   (3 separate runs)
$p = New-ParaAfter $rng "- Xem các tour "
$p = Add-Run $p "sắp có"
$p = Add-Run $p "."

# 3b) "- Liên hệ với công ty"  (the _GoBack bookmark will be re-attached
#      to the end of this paragraph afterwards)
$p = New-ParaAfter $p "- Liên hệ với công ty"

# 3c) "- " / "Xem thông tin tour mình đã đặt."   (2 separate runs)
$p = New-ParaAfter $p "- "
$p = Add-Run $p "Xem thông tin tour mình đã đặt."

# 3d) "- Xem thông tin tour mình đã đi."
$p = New-ParaAfter $p "- Xem thông tin tour mình đã đi."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Host "$i => [$($pp.Range.Text)]"
}
